$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.534.19'
$ws.Range("E2").Value = '  +1.72%  '

$ws.Range("D3").Value = '2.604.99'
$ws.Range("E3").Value = '  +0.68%  '

$ws.Range("E4").Value = '  +0.01%  '

$c = $ws.Range("D5")
$c.Value = "'573.47"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.12%  '

$c = $ws.Range("D6")
$c.Value = "'142.97"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").Value = '  +0.51%  '

$ws.Range("D9").Value = '2.629.02'
$ws.Range("E9").Value = '  +1.25%  '

$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("E11").Value = '  +1.15%  '

$c = $ws.Range("D12")
$c.Value = "'0.152"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.88%  '

$ws.Range("E13").Value = '  +2.23%  '

$ws.Range("D14").Value = '3.072.82'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '60.534.27'
$ws.Range("E15").Value = '  +1.74%  '

$ws.Range("E16").Value = '  -0.30%  '

$ws.Range("D18").Value = '2.618.65'
$ws.Range("E18").Value = '  +0.94%  '

$c = $ws.Range("D19")
$c.Value = "'11.34"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +9.02%  '

$ws.Range("E20").Value = '  +1.70%  '

$c = $ws.Range("D21")
$c.Value = "'347.55"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.51%  '

$c = $ws.Range("D22")
$c.Value = "'6.95"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +6.64%  '

$ws.Range("E23").Value = '  -0.39%  '

$c = $ws.Range("D24")
$c.Value = "'0.528"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +12.00%  '

$c = $ws.Range("D25")
$c.Value = "'63.25"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  +0.01%  '

$c = $ws.Range("D27")
$c.Value = "'0.160"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("E28").Value = '  +3.94%  '

$ws.Range("D29").Value = '0.0₃0790'
$ws.Range("E29").Value = '  +1.24%  '

$ws.Range("E30").Value = '  +10.29%  '

$c = $ws.Range("D31")
$c.Value = "'6.38"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.48%  '

$ws.Range("E32").Value = '  -0.10%  '

$c = $ws.Range("D33")
$c.Value = "'161.91"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.36%  '

$ws.Range("E34").Value = '  +2.25%  '

$ws.Range("E35").Value = '  +4.17%  '

$c = $ws.Range("D36")
$c.Value = "'0.980"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +8.67%  '

$c = $ws.Range("D37")
$c.Value = "'1.23"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +4.58%  '

$ws.Range("E38").Value = '  +7.65%  '

$ws.Range("E39").Value = '  +1.11%  '

$c = $ws.Range("D40")
$c.Value = "'3.84"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.40%  '

$c = $ws.Range("D41")
$c.Value = "'0.842"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.11%  '

$c = $ws.Range("D42")
$c.Value = "'295.21"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '

$c = $ws.Range("D43")
$c.Value = "'137.33"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.15%  '

$ws.Range("E44").Value = '  -0.27%  '

$c = $ws.Range("D45")
$c.Value = "'0.0985"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.90%  '

$c = $ws.Range("D46")
$c.Value = "'0.608"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '

$ws.Range("E47").Value = '  +3.31%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D48")
$c.Value = "'4.96"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +9.20%  '

$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D49")
$c.Value = "'0.0544"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.19%  '

$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("E51").Value = '  +0.52%  '
